# Request Report.xlsx - "Reports" sheet update
# - FindPage: swap the Request ID / Report Timestamp header columns
# - ReceivePage: replace the dummy row 2 with the real "found battery" report
# - add row 3 for the "received battery" report (PicturePage still pending)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: Request ID now comes before Report Timestamp ---
$ws.Range("A1").Value = "Request ID"
$ws.Range("B1").Value = "Report Timestamp"
$ws.Range("C1").Value = "Report"

# --- Row 2: first report (battery found) ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "2025-04-25 10:16:26"
$ws.Range("C2").Value = "John Smith found battery 2. Now John Smith is Frustrated"

# --- Row 3 (new): second report (battery received) ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "2025-04-25 10:17:01"
$ws.Range("C3").Value = "John Smith received battery 2 from Battery New. Now John Smith is Tired."

# --- Column widths: Request ID / Report Timestamp / Report ---
# (values chosen so the stored OOXML <col width> lands on 14.4 / 25.2 / 88.8,
#  compensating for this host's char-width -> stored-width rounding)
$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Columns.Item(2).ColumnWidth = 24.333333333333332
$ws.Columns.Item(3).ColumnWidth = 88.0

# --- Selection moves to C9 ---
$null = $ws.Range("C9").Select()
